# Daily attendance processing - 2025-10-14 18:29:16
# Re-order "Recorded By" email lists, tweak the summary metric counts,
# and flip session row 40 (Year3/C2/PHARMACOLOGY #3) from "Not Recorded"
# (pink) back to "Pending" (yellow) with a corrected date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recorded-by email list re-orderings ---------------------------------

$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"

$ws.Range("G12").Value = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

$ws.Range("G19").Value = "ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

$ws.Range("G25").Value = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

$ws.Range("G26").Value = "hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"

$ws.Range("G34").Value = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"

$ws.Range("G41").Value = "Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

$ws.Range("G42").Value = "ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

# --- Summary metric table tweaks ------------------------------------------

$ws.Range("L7").Value = 4
$ws.Range("L8").Value = 29

$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 14

# --- Row 40 content fix: status text + date --------------------------------
# The date string is textual (not a real date), so force literal text entry
# with a leading apostrophe to stop Excel from auto-converting it to a
# date serial number, exactly as a user re-typing it would do.

$ws.Range("E40").Value = "'10/11/2025"
$ws.Range("I40").Value = "Pending"

# --- Row 40 style fix: was "Not Recorded" (pink, style of row 3/27/etc.),
# now "Pending" (yellow, same look as rows 7/8/9/16/17/18...). Copy the
# cell formatting from row 7 (a known "Pending" row) across A:I so the
# whole row matches exactly, including the just-edited E40/I40 cells.

$ws.Range("A7:I7").Copy()
$ws.Range("A40:I40").PasteSpecial(-4122)
